# "Updated for first release" - refresh Jacobi benchmark numbers, switch the
# active sheet/selection back to Jacobi, convert the F2:F10 runtime column to
# a shared formula, and drop the legend from the first (single-series)
# Runtime chart.

$wb = $excel.ActiveWorkbook

$jacobi = $wb.Worksheets.Item("Jacobi")

# --- Jacobi!F2:F10 -> shared formula (same formulas, just re-entered as a
#     block so the writer collapses them into t="shared") ---------------
$jacobi.Range("F2:F10").Formula = "=B2"

# --- Jacobi row 11 (the "Kernels" row) got new measured numbers --------
$jacobi.Range("B11").Value = 74.480391999999995
$jacobi.Range("C11").Value = 22.91516
$jacobi.Range("D11").Value = 11.467739999999999
$jacobi.Range("F11").Formula = "=2.160203+1.232692+0.221674"
# E11 (=C11+D11) and G11 (=B11-(E11+F11)) recompute on their own.

# --- Chart 1 ("Runtime") loses its legend --------------------------------
$chart1 = $jacobi.ChartObjects(1).Chart
$chart1.HasLegend = $false

# --- Active sheet / selection moves back to Jacobi!B12 ------------------
# (MatVec's own selection, E30, is left untouched - it just stops being the
# active tab once Jacobi is activated below.)
[void]$jacobi.Activate()
[void]$jacobi.Range("B12").Select()
